# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
#
# The source data rows for three fixtures got re-matched to the correct
# match ids / odds, which - for this sheet - manifests as the full data
# (columns B:AB) of each of the following row pairs being swapped with
# each other, while the row's own sequence number in column A stays put:
#   rows 104 <-> 105
#   rows 112 <-> 113
#   rows 124 <-> 125
# Additionally the closing Asian-handicap odds for the still-unplayed
# fixture on row 171 were updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, [int]$rowA, [int]$rowB) {
    $rangeA = $sheet.Range("B$rowA`:AB$rowA")
    $rangeB = $sheet.Range("B$rowB`:AB$rowB")
    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()
    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

Swap-Rows $ws 104 105
Swap-Rows $ws 112 113
Swap-Rows $ws 124 125

# Row 171 (future fixture, id 169): closing AH odds update
$ws.Range("Q171").Value = 1.89
$ws.Range("R171").Value = 2.01
